$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("association")
$ws.Activate()

# Match A7's style for A8 (copy format from A7, then set value)
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C7").Value = "1"
$ws.Range("B7").Value = "no-rsid-val"
$ws.Range("F8").Value = "T"
$ws.Range("C8").Value = "1E-5"
$ws.Range("B8").Value = "no-rsid-val"
$ws.Range("A8").Value = "white"

$ws.Range("D8").Select()
